# Update cryptos list (price + 1h volume change) - GitHub Actions data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store literal text (e.g. "1.019", "27.916.52")
# rather than numbers, so force text formatting before writing the new values,
# then restore the original (default) style once done.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '27.916.52'
$ws.Range("E2").Value = '  +0.74%  '

# Row 3
$ws.Range("D3").Value = '1.877.43'
$ws.Range("E3").Value = '  -0.05%  '

# Row 4
$ws.Range("E4").Value = '  +1.64%  '

# Row 5
$ws.Range("D5").Value = '334.76'
$ws.Range("E5").Value = '  +1.11%  '

# Row 6
$ws.Range("D6").Value = '1.019'
$ws.Range("E6").Value = '  +1.66%  '

# Row 7
$ws.Range("D7").Value = '0.4691'
$ws.Range("E7").Value = '  -0.29%  '

# Row 8
$ws.Range("D8").Value = '0.3912'
$ws.Range("E8").Value = '  -1.27%  '

# Row 9
$ws.Range("D9").Value = '46.59'
$ws.Range("E9").Value = '  -2.98%  '

# Row 10
$ws.Range("D10").Value = '0.07938'
$ws.Range("E10").Value = '  -1.15%  '

# Row 11
$ws.Range("D11").Value = '1.005'
$ws.Range("E11").Value = '  -1.84%  '

# Row 12
$ws.Range("D12").Value = '21.56'
$ws.Range("E12").Value = '  -1.39%  '

# Row 13
$ws.Range("D13").Value = '1.892.15'
$ws.Range("E13").Value = '  -0.72%  '

# Row 14
$ws.Range("D14").Value = '5.932'

# Row 15
$ws.Range("D15").Value = '7.091'
$ws.Range("E15").Value = '  -1.15%  '

# Row 16
$ws.Range("D16").Value = '1.021'
$ws.Range("E16").Value = '  +1.55%  '

# Row 17
$ws.Range("D17").Value = '0.06780'
$ws.Range("E17").Value = '  +2.50%  '

# Row 18
$ws.Range("D18").Value = '87.32'
$ws.Range("E18").Value = '  +0.19%  '

# Row 19
$ws.Range("D19").Value = '0.00001043'
$ws.Range("E19").Value = '  -0.07%  '

# Row 20
$ws.Range("D20").Value = '16.99'
$ws.Range("E20").Value = '  -1.77%  '

# Row 21
$ws.Range("D21").Value = '1.017'
$ws.Range("E21").Value = '  +1.49%  '

# Row 22
$ws.Range("D22").Value = '27.913.86'
$ws.Range("E22").Value = '  +0.69%  '

# Row 23
$ws.Range("D23").Value = '5.461'
$ws.Range("E23").Value = '  -0.79%  '

# Row 24
$ws.Range("D24").Value = '10.93'
$ws.Range("E24").Value = '  -1.03%  '

# Row 25
$ws.Range("D25").Value = '2.364'
$ws.Range("E25").Value = '  +2.91%  '

# Row 26
$ws.Range("D26").Value = '2.140.08'
$ws.Range("E26").Value = '  +0.60%  '

# Row 27
$ws.Range("D27").Value = '159.95'
$ws.Range("E27").Value = '  +2.18%  '

# Row 28
$ws.Range("D28").Value = '19.85'
$ws.Range("E28").Value = '  -1.95%  '

# Row 29
$ws.Range("D29").Value = '2.069'
$ws.Range("E29").Value = '  -1.21%  '

# Row 30
$ws.Range("D30").Value = '5.443'
$ws.Range("E30").Value = '  -2.64%  '

# Row 31
$ws.Range("D31").Value = '120.83'
$ws.Range("E31").Value = '  -1.47%  '

# Row 32
$ws.Range("D32").Value = '0.09503'
$ws.Range("E32").Value = '  -0.52%  '

# Row 33
$ws.Range("D33").Value = '0.9564'
$ws.Range("E33").Value = '  -1.57%  '

# Row 34
$ws.Range("D34").Value = '3.667'
$ws.Range("E34").Value = '  +1.08%  '

# Row 35
$ws.Range("D35").Value = '5.297'
$ws.Range("E35").Value = '  -0.16%  '

# Row 36
$ws.Range("D36").Value = '1.341'
$ws.Range("E36").Value = '  -7.72%  '

# Row 37
$ws.Range("D37").Value = '0.06108'
$ws.Range("E37").Value = '  -0.14%  '

# Row 38
$ws.Range("D38").Value = '0.02237'
$ws.Range("E38").Value = '  -1.18%  '

# Row 39
$ws.Range("D39").Value = '1.200'
$ws.Range("E39").Value = '  -2.78%  '

# Row 40
$ws.Range("D40").Value = '1.018'
$ws.Range("E40").Value = '  +1.61%  '

# Row 41
$ws.Range("D41").Value = '8.099'
$ws.Range("E41").Value = '  -0.87%  '

# Row 42
$ws.Range("D42").Value = '0.5897'
$ws.Range("E42").Value = '  -1.78%  '

# Row 43
$ws.Range("D43").Value = '0.1885'
$ws.Range("E43").Value = '  -0.81%  '

# Row 44
$ws.Range("D44").Value = '10.14'
$ws.Range("E44").Value = '  -0.88%  '

# Row 45
$ws.Range("D45").Value = '1.271'
$ws.Range("E45").Value = '  +1.67%  '

# Row 46
$ws.Range("D46").Value = '0.5634'
$ws.Range("E46").Value = '  -1.05%  '

# Row 47
$ws.Range("D47").Value = '12.07'
$ws.Range("E47").Value = '  -1.73%  '

# Row 48
$ws.Range("D48").Value = '3.410'
$ws.Range("E48").Value = '  +0.36%  '

# Row 49
$ws.Range("D49").Value = '1.915'
$ws.Range("E49").Value = '  -0.98%  '

# Row 50
$ws.Range("D50").Value = '0.06867'
$ws.Range("E50").Value = '  +0.67%  '

# Row 51
$ws.Range("D51").Value = '113.65'
$ws.Range("E51").Value = '  +1.79%  '

# Restore default (General) style now that values are written as text
$ws.Range("D2:E51").Style = "Normal"
